{"js": "// Replace \"passphrase\" with \"letters\" wherever it appears (the k3 parameter\n// descriptions), matching the author's intent in the commit message:\n// \"Changed wording, 'passphrase' changed to 'letters' in all spots.\"\nconst body = context.document.body;\nconst results = body.search(\"passphrase\", { matchCase: false, matchWholeWord: true });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  const target = results.items[i];\n  target.insertText(\"letters\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Re-locate the freshly inserted word(s) and nudge formatting on/off so the\n// run boundaries around the replaced word stay distinct (matching how Word\n// itself leaves the surrounding text split into separate runs instead of\n// silently re-merging them back into one run).\nconst replaced = body.search(\"letters\", { matchCase: false, matchWholeWord: true });\nreplaced.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < replaced.items.length; i++) {\n  const r = replaced.items[i];\n  r.font.load(\"bold\");\n  await context.sync();\n  const original = r.font.bold;\n  r.font.bold = !original;\n  await context.sync();\n  r.font.bold = original;\n  await context.sync();\n}\n", "ps1": "# Replace \"passphrase\" with \"letters\" wherever it appears (the k3 parameter\n# descriptions), matching the commit message:\n# \"Changed wording, 'passphrase' changed to 'letters' in all spots.\"\n$d = $word.ActiveDocument\n$searchRange = $d.Content\n\n$found = $true\nwhile ($found) {\n    $find = $searchRange.Find\n    $find.ClearFormatting()\n    $find.Text = \"passphrase\"\n    $find.MatchWholeWord = $true\n    $find.MatchCase = $false\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = \"letters\"\n    $found = $find.Execute([ref]\"passphrase\", [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]\"letters\", [ref]1)\n\n    if ($found) {\n        # After Execute, $searchRange now spans the just-replaced \"letters\"\n        # text. Nudge a character property on and back off so Word keeps\n        # the surrounding sentence split into its own runs around the\n        # replaced word, instead of silently re-merging everything back\n        # into a single run.\n        $originalBold = $searchRange.Bold\n        $searchRange.Bold = 1 - $originalBold\n        $searchRange.Bold = $originalBold\n\n        # Continue searching after this occurrence.\n        $searchRange.Start = $searchRange.End\n        $searchRange.End = $d.Content.End\n    }\n}\n"}
